$wb = $excel.ActiveWorkbook

# --- Overview sheet: update per-locale status for the two rows that moved
#     from "Ready for handoff" to "In Translation" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"

# --- zh-cn sheet: same two rows (identified by source file), Status column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("C4").Value = "In Translation"

# --- de-de sheet: same two rows, Status column C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("C4").Value = "In Translation"
